$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D ("password" header + "employee123" value) and shift cells left.
$ws.Range("D1:D2").EntireColumn.Delete()

# Update the active cell selection to D1 (per diff).
$ws.Range("D1").Select()
